$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "PERALTA REYES MARY CRUZ",
    "ESPINOZA GUZMAN MAYRA LOURDES",
    "RUIZ CARRASCO HILLARY SAMANTHA",
    "CHUNGA DE LA CRUZ ROSA LILIANA",
    "FIESTAS PERICHE VIVIANA LISSETH",
    "ROSILLO ALBERCA ROXANA",
    "PAIVA PINDAY ALICIA",
    "PRADO ACARO VANESSA PAOLA",
    "GIRON SILUPU JUAN FRANCISCO",
    "PALMA CARMENES DE MENA MERCEDES EVERJISTA",
    "JIMENEZ GUERRERO JUAN RICARDO",
    "PAIVA GARCIA DANIELA MILEYDI"
)

$values = @(180, 121, 119, 118, 116, 94, 77, 73, 70, 68, 64, 53)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
